# New PO forecast model
# Updates the three worksheets (Weekly Quantity, Monthly Trend, PO Forecast)
# with the latest pulled PO data / forecast recompute.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Weekly Quantity": append the latest observed week
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Range("A21").NumberFormat = $ws1.Range("A20").NumberFormat
$ws1.Range("A21").Value = 45676.99999999999
$ws1.Range("B21").Value = 1

# ---------------------------------------------------------------------------
# Sheet "Monthly Trend": append the latest observed month
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("A9").NumberFormat = $ws2.Range("A8").NumberFormat
$ws2.Range("A9").Value = 45688.99999999999
$ws2.Range("B9").Value = 1

# ---------------------------------------------------------------------------
# Sheet "PO Forecast": recomputed forecast values, and extended horizon
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("PO Forecast")

$ws3.Range("B2").Value = 26
$ws3.Range("B3").Value = 26
$ws3.Range("B4").Value = 26
$ws3.Range("B5").Value = 26
$ws3.Range("B6").Value = 26
$ws3.Range("B7").Value = 26
$ws3.Range("B8").Value = 26
$ws3.Range("B9").Value = 26
$ws3.Range("B10").Value = 26
$ws3.Range("B12").Value = 26
$ws3.Range("B13").Value = 26
$ws3.Range("B14").Value = 26
$ws3.Range("B15").Value = 26
$ws3.Range("B16").Value = 25
$ws3.Range("B17").Value = 25
$ws3.Range("B18").Value = 25
$ws3.Range("B19").Value = 25
$ws3.Range("B20").Value = 25

$ws3.Range("A21").Value = 45676.99999999999
$ws3.Range("B21").Value = 25
$ws3.Range("A22").Value = 45683.99999999999
$ws3.Range("B22").Value = 25
$ws3.Range("A23").Value = 45690.99999999999
$ws3.Range("B23").Value = 25
$ws3.Range("A24").Value = 45697.99999999999
$ws3.Range("B24").Value = 25
$ws3.Range("A25").Value = 45704.99999999999
$ws3.Range("B25").Value = 25
$ws3.Range("A26").Value = 45711.99999999999
$ws3.Range("B26").Value = 25
$ws3.Range("A27").Value = 45718.99999999999
$ws3.Range("B27").Value = 25
$ws3.Range("A28").Value = 45725.99999999999
$ws3.Range("B28").Value = 25

# New row 29 - extend the forecast horizon one more week
$ws3.Range("A29").NumberFormat = $ws3.Range("A28").NumberFormat
$ws3.Range("A29").Value = 45732.99999999999
$ws3.Range("B29").Value = 25
